# Generate Report for Handback
# -----------------------------------------------------------------------
# This script mirrors the localization-status report being regenerated
# after a handback: the "zh-cn" and "de-de" detail sheets get their
# "Latest Target File" / "Latest Handback File" / "Latest Handback
# DateTime" columns populated (they were empty / placeholder before),
# and the overview "Status" column flips from "Ready for handoff" to
# "Handed back: in sync with en-US". Columns are widened to fit the
# newly-populated, longer text.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$repoBlobBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ed32bfdc2ea3787fea8841288ae3ff020df28854/e2e/"
$mdRow2 = "de25e7ca-e1a6-4890-ad3d-9c70ce26b4eb.md"
$mdRow2Url = $repoBlobBase + $mdRow2

$statusHandedBack = "Handed back: in sync with en-US"

# -----------------------------------------------------------------------
# Overview sheet: Status columns (zh-cn / de-de) now read "Handed back"
# -----------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusHandedBack
$overview.Range("F2").Value = $statusHandedBack
$overview.Range("E3").Value = $statusHandedBack
$overview.Range("F3").Value = $statusHandedBack

$overview.Columns.Item(5).ColumnWidth = 29.15
$overview.Columns.Item(6).ColumnWidth = 29.15

# -----------------------------------------------------------------------
# zh-cn detail sheet
# -----------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

# "Status" column also flips to the handed-back message
$zhcn.Range("C2").Value = $statusHandedBack
$zhcn.Range("C3").Value = $statusHandedBack

# Row 2 (de25e7ca....md)
$zhcn.Range("I2").Value = $mdRow2
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $mdRow2Url, "", "", $mdRow2)
$zhcn.Range("J2").Value = "de25e7ca-e1a6-4890-ad3d-9c70ce26b4eb.551d55050f20a10c8fd718c21c869ebd7bafc806.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-24 21:05:34"

# Row 3 (ffff98a7d8fa-....md) also got handed back for the same target file
$zhcn.Range("I3").Value = $mdRow2
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $mdRow2Url, "", "", $mdRow2)
$zhcn.Range("J3").Value = "de25e7ca-e1a6-4890-ad3d-9c70ce26b4eb.551d55050f20a10c8fd718c21c869ebd7bafc806.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-24 21:05:34"

$zhcn.Columns.Item(3).ColumnWidth = 29.15
$zhcn.Columns.Item(9).ColumnWidth = 39.17
$zhcn.Columns.Item(10).ColumnWidth = 39.17

# -----------------------------------------------------------------------
# de-de detail sheet
# -----------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = $statusHandedBack
$dede.Range("C3").Value = $statusHandedBack

# Row 2
$dede.Range("I2").Value = $mdRow2
$dede.Hyperlinks.Add($dede.Range("I2"), $mdRow2Url, "", "", $mdRow2)
$dede.Range("J2").Value = "de25e7ca-e1a6-4890-ad3d-9c70ce26b4eb.551d55050f20a10c8fd718c21c869ebd7bafc806.de-de.xlf"
$dede.Range("K2").Value = "2016-08-24 21:05:43"

# Row 3
$dede.Range("I3").Value = $mdRow2
$dede.Hyperlinks.Add($dede.Range("I3"), $mdRow2Url, "", "", $mdRow2)
$dede.Range("J3").Value = "de25e7ca-e1a6-4890-ad3d-9c70ce26b4eb.551d55050f20a10c8fd718c21c869ebd7bafc806.de-de.xlf"
$dede.Range("K3").Value = "2016-08-24 21:05:43"

$dede.Columns.Item(3).ColumnWidth = 29.15
$dede.Columns.Item(9).ColumnWidth = 39.17
$dede.Columns.Item(10).ColumnWidth = 39.17
